$d = $word.ActiveDocument

$insertPos = 3844
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter(" ")

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Step 1: add bookmark around space to force split (creates bookmarkStart...space...bookmarkEnd)
$bmRange = $d.Range(3844, 3845)
$d.Bookmarks.Add("_GoBack", $bmRange)
$b = $d.Bookmarks.Item("_GoBack")
Write-Host "After step1 GoBack range: $($b.Range.Start)-$($b.Range.End)"

# Step 2: delete and re-add collapsed at end (3845) now that a run boundary exists there
$d.Bookmarks.Item("_GoBack").Delete()
$bmRange2 = $d.Range(3845, 3845)
Write-Host "bmRange2: $($bmRange2.Start)-$($bmRange2.End)"
$d.Bookmarks.Add("_GoBack", $bmRange2)
$b2 = $d.Bookmarks.Item("_GoBack")
Write-Host "After step2 GoBack range: $($b2.Range.Start)-$($b2.Range.End)"
